$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 4 Remarks as "Done"
$ws.Range("J4").Value = "Done"

# Fill in row 6 with the new time-sheet entry
$ws.Range("B6").Value = "17-12-2018"
$ws.Range("C6").Value = "D.Venkatesh"
$ws.Range("D6").Value = "Twitter Data Analytics & Content preparation for AI & ML, Data Science & ML , Python & Working On Telagana Whatsapp Project"
$ws.Range("E6").Value = "Pushing the Twitter data in Database using MSSQL Server & Content Preparation & Coordinating with Hyderabad Team."
$ws.Range("H6").Value = "9:00"
$ws.Range("I6").Value = "5:00"
$ws.Range("J6").Value = "1)Content Preparation Done.                                                2) Waiting for whatsapp Numbers & Whatsapp Work will Starts from Tommrrow"

# Grow row 6 to fit the new wrapped remarks text
$ws.Rows.Item(6).RowHeight = 72

# Move the active selection to C8, matching the author's last edit position
$ws.Range("C8").Select()
